# Apply the "handles float input without breaking stuff" marksheet fix.
# The quiz result sheet previously tracked three separate attempt columns
# (A/B, D/E, G/H) with placeholder/blank "Student Ans" data and bogus
# summary numbers. This edit:
#   1. Fixes the summary row numbers (Right/Wrong/Not-Attempt/Max/Total).
#   2. Fixes the marking value & Wrong cell numeric type (was text "-1").
#   3. Removes the third attempt block (columns G:H) entirely.
#   4. Removes the stray D:E "not attempted" placeholders for rows 19-40
#      (leaving just the single remaining attempt in columns A:B).
#   5. Fills in the actual "Student Ans" values (column A, and for a few
#      rows column D) with the right grading style (correct/incorrect).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: summary rows 10-12 -------------------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Style = "incorrectStyle"
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "57/112"

# --- 3: drop the third "Student Ans / Correct Ans" attempt (columns G:H) ------
$thirdAttemptCells = @(
    "G15","H15",
    "G16","H16",
    "G17","H17",
    "G18","H18",
    "G19","H19",
    "G20","H20",
    "G21","H21"
)
foreach ($addr in $thirdAttemptCells) {
    $ws.Range($addr).Clear()
}

# --- 4: drop the stray second-attempt placeholders for rows 19-40 -------------
$secondAttemptCells = @(
    "D19","E19",
    "D20","E20",
    "D21","E21",
    "D22","E22",
    "D23","E23",
    "D24","E24",
    "D25","E25",
    "D26","E26",
    "D27","E27",
    "D28","E28",
    "D29","E29",
    "D30","E30",
    "D31","E31",
    "D32","E32",
    "D33","E33",
    "D34","E34",
    "D35","E35",
    "D36","E36",
    "D37","E37",
    "D38","E38",
    "D39","E39",
    "D40","E40"
)
foreach ($addr in $secondAttemptCells) {
    $ws.Range($addr).Clear()
}

# --- 5: fill in the real "Student Ans" values ----------------------------------
# Remaining second-attempt cells (rows 16-18) that keep their "Student Ans" value
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

$ws.Range("D17").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"

$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"

# Student answers for the primary attempt (column A), rows 19-40.
# normalStyle (blank) = not attempted -> leave as-is; correctStyle/incorrectStyle
# = attempted, graded right/wrong.
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"

$ws.Range("A22").Style = "incorrectStyle"
$ws.Range("A22").Value = "Option A"

$ws.Range("A23").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"

$ws.Range("A24").Style = "correctStyle"
$ws.Range("A24").Value = "Option A"

$ws.Range("A26").Style = "correctStyle"
$ws.Range("A26").Value = "Option C"

$ws.Range("A27").Style = "incorrectStyle"
$ws.Range("A27").Value = "Option C"

$ws.Range("A29").Style = "correctStyle"
$ws.Range("A29").Value = "Option D"

$ws.Range("A31").Style = "correctStyle"
$ws.Range("A31").Value = "Option D"

$ws.Range("A32").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"

$ws.Range("A34").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"

$ws.Range("A36").Style = "correctStyle"
$ws.Range("A36").Value = "Option A"

$ws.Range("A37").Style = "correctStyle"
$ws.Range("A37").Value = "Option A"

$ws.Range("A38").Style = "incorrectStyle"
$ws.Range("A38").Value = "Option C"

$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
